$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet/tab name to reflect the new "through" date
$ws.Name = "Through 2022-09-21"

# Update the row label for September to reflect the new "through" date
$ws.Range("A10").Value = "September (through 09-21)"

# Update September row (row 10) values for each year column (B..I = 2015..2022)
$ws.Range("B10").Value = 23
$ws.Range("C10").Value = 34
$ws.Range("D10").Value = 49
$ws.Range("E10").Value = 40
$ws.Range("F10").Value = 50
$ws.Range("G10").Value = 82
$ws.Range("H10").Value = 126
$ws.Range("I10").Value = 102

# Update Total row (row 11) values for each year column (B..I = 2015..2022)
$ws.Range("B11").Value = 217
$ws.Range("C11").Value = 415
$ws.Range("D11").Value = 600
$ws.Range("E11").Value = 530
$ws.Range("F11").Value = 399
$ws.Range("G11").Value = 866
$ws.Range("H11").Value = 1196
$ws.Range("I11").Value = 1237
